$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 38, shifting rows 38-41 down to 39-42.
$ws.Rows.Item(38).Insert()

# Update E37 value (this affects the shared formulas F37/G37 already present).
$ws.Range("E37").Value = 0.72916666666666663

# Fill in new row 38 data.
$ws.Range("A38").Value = 2014
$ws.Range("B38").Value = 3
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 0.74305555555555547
$ws.Range("E38").Value = 0.77083333333333337
$ws.Range("F38").Formula = "=(E38-D38)*24*60"
$ws.Range("G38").Formula = "=F38/60"

$ws.Range("F40").Formula = "=SUM(F2:F39)"
$ws.Range("F41").Formula = "=F40/60"
$ws.Range("F42").Formula = "=F41/38.5"

$ws.Range("F38").Select()
